# Updates crypto price/volume data to the latest scraped values.
# Some "Price" column entries are plain decimal-looking strings (e.g. "245.84")
# that must remain TEXT (matching the sheet's existing inline-string convention)
# rather than being auto-converted to numbers by Excel's type inference.
# For those, force the cell to Text format before writing, then restore the
# default "Normal" style so no stray formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.375.56'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '1.909.61'
$ws.Range("E3").Value = '  +2.77%  '
$ws.Range("E4").Value = '  -0.48%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.84'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.661'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.31%  '
$ws.Range("E7").Value = '  -0.49%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.36'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.85%  '
$ws.Range("E9").Value = '  +6.50%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '52.80'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +12.54%  '
$ws.Range("E11").Value = '  +3.84%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0992'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.36%  '
$ws.Range("D13").Value = '2.186.86'
$ws.Range("E13").Value = '  +2.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '12.13'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.67%  '
$ws.Range("E15").Value = '  +3.69%  '
$ws.Range("D16").Value = '1.898.05'
$ws.Range("E16").Value = '  +3.38%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.87'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.24%  '
$ws.Range("D18").Value = '35.350.40'
$ws.Range("E18").Value = '  +0.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.71'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.10%  '
$ws.Range("D20").Value = '0.0₃0824'
$ws.Range("E20").Value = '  +3.65%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '239.83'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("E22").Value = '  +2.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.91%  '
$ws.Range("E24").Value = '  -0.46%  '
$ws.Range("E25").Value = '  +1.37%  '
$ws.Range("E26").Value = '  +23.96%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '169.91'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.34%  '
$ws.Range("E28").Value = '  +5.84%  '
$ws.Range("E29").Value = '  +4.73%  '
$ws.Range("E30").Value = '  +2.43%  '
$ws.Range("E31").Value = '  +4.01%  '
$ws.Range("E32").Value = '  +1.37%  '
$ws.Range("B33").Value = 'BinanceUSD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.01'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.42%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.939'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +15.73%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.12'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.71%  '
$ws.Range("E36").Value = '  -4.03%  '
$ws.Range("E37").Value = '  +0.62%  '
$ws.Range("E38").Value = '  +0.79%  '
$ws.Range("E39").Value = '  +1.65%  '
$ws.Range("E40").Value = '  +3.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '16.26'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.54%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0635'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '90.08'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("D44").Value = '1.339.77'
$ws.Range("E44").Value = '  -0.46%  '
$ws.Range("E45").Value = '  +3.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '47.96'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +38.27%  '
$ws.Range("E47").Value = '  +1.72%  '
$ws.Range("E48").Value = '  -0.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.56'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.58%  '
$ws.Range("D50").Value = '2.093.46'
$ws.Range("E50").Value = '  +2.39%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0707'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.02%  '
